$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "Resolving-Mac" sending-cluster block (old rows 14-17)
$ws.Range("A14:T17").EntireRow.Delete()

# Re-run NATMI computation results with updated TPM values (rows 2-13)
$data = @(
  @("ECs", "Vtn", "Plaur", "ECs", 3, 1, 7.844453333333334, 23.53336, 0.1489290605659587, 0.1489290605659588, 3, 1, 15.03463666666667, 45.10391, 0.2402934356091235, 0.2402934356091235, 117.9385057152889, 1061.4465514376, 0.03578667562543346, 0.03578667562543347),
  @("ECs", "Vtn", "Plaur", "FAPs", 3, 1, 7.844453333333334, 23.53336, 0.1489290605659587, 0.1489290605659588, 3, 1, 4.827410666666666, 14.482232, 0.07715484716443403, 0.07715484716443403, 37.86839769550222, 340.81557925952, 0.01149059890630929, 0.01149059890630929),
  @("ECs", "Vtn", "Plaur", "MuSCs", 3, 1, 7.844453333333334, 23.53336, 0.1489290605659587, 0.1489290605659588, 3, 1, 6.211932333333333, 18.635797, 0.09928318157880762, 0.09928318157880764, 48.72921329865778, 438.56291968792, 0.01478615096253132, 0.01478615096253132),
  @("ECs", "Vtn", "Plaur", "Resolving-Mac", 3, 1, 7.844453333333334, 23.53336, 0.1489290605659587, 0.1489290605659588, 3, 1, 36.49384133333334, 109.481524, 0.5832685356476348, 0.5832685356476348, 286.2742352934045, 2576.46811764064, 0.08686563507168465, 0.08686563507168468),
  @("FAPs", "Vtn", "Plaur", "ECs", 3, 1, 20.35396833333334, 61.06190500000001, 0.3864255740794268, 0.3864255740794268, 3, 1, 15.03463666666667, 45.10391, 0.2402934356091235, 0.2402934356091235, 306.0145186165056, 2754.13066754855, 0.09285552880277333, 0.09285552880277334),
  @("FAPs", "Vtn", "Plaur", "FAPs", 3, 1, 20.35396833333334, 61.06190500000001, 0.3864255740794268, 0.3864255740794268, 3, 1, 4.827410666666666, 14.482232, 0.07715484716443403, 0.07715484716443403, 98.2569638413289, 884.3126745719601, 0.02981460610852685, 0.02981460610852685),
  @("FAPs", "Vtn", "Plaur", "MuSCs", 3, 1, 20.35396833333334, 61.06190500000001, 0.3864255740794268, 0.3864255740794268, 3, 1, 6.211932333333333, 18.635797, 0.09928318157880762, 0.09928318157880764, 126.4374740014761, 1137.937266013285, 0.03836556043802271, 0.03836556043802272),
  @("FAPs", "Vtn", "Plaur", "Resolving-Mac", 3, 1, 20.35396833333334, 61.06190500000001, 0.3864255740794268, 0.3864255740794268, 3, 1, 36.49384133333334, 109.481524, 0.5832685356476348, 0.5832685356476348, 742.794490860358, 6685.150417743222, 0.2253898787301039, 0.2253898787301039),
  @("MuSCs", "Vtn", "Plaur", "ECs", 3, 1, 24.47399366666667, 73.421981, 0.4646453653546145, 0.4646453653546145, 3, 1, 15.03463666666667, 45.10391, 0.2402934356091235, 0.2402934356091235, 367.9576025606344, 3311.61842304571, 0.1116512311809167, 0.1116512311809167),
  @("MuSCs", "Vtn", "Plaur", "FAPs", 3, 1, 24.47399366666667, 73.421981, 0.4646453653546145, 0.4646453653546145, 3, 1, 4.827410666666666, 14.482232, 0.07715484716443403, 0.07715484716443403, 118.1460180823991, 1063.314162741592, 0.03584964214959789, 0.0358496421495979),
  @("MuSCs", "Vtn", "Plaur", "MuSCs", 3, 1, 24.47399366666667, 73.421981, 0.4646453653546145, 0.4646453653546145, 3, 1, 6.211932333333333, 18.635797, 0.09928318157880762, 0.09928318157880764, 152.0307925837619, 1368.277133253857, 0.0461314701782536, 0.04613147017825361),
  @("MuSCs", "Vtn", "Plaur", "Resolving-Mac", 3, 1, 24.47399366666667, 73.421981, 0.4646453653546145, 0.4646453653546145, 3, 1, 36.49384133333334, 109.481524, 0.5832685356476348, 0.5832685356476348, 893.1500416643383, 8038.350374979045, 0.2710130218458462, 0.2710130218458463),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $col = $j + 1
        $ws.Cells.Item($excelRow, $col).Value = $row[$j]
    }
}
